# EPBDS-7993 Use [] to escape special symbols
#
# Rows 22-41 of column C/D/E on sheet "like" hold test data (String, Pattern, Result)
# for the like() function. The three "escaped @" pattern rows get rewritten to use the
# new [@] bracket-escape syntax, and the test rows get reordered accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for rows 22..41: Column C (string), Column D (pattern), Column E (bool result)
$data = @(
    @{ Row = 22; C = "foo.bar@gmail.";     D = "?+[@]?+.?+";       E = $false },
    @{ Row = 23; C = "foo.bar@gmailcom";   D = "?+[@]?+.?+";       E = $false },
    @{ Row = 24; C = "foo.bar@gmail.com";  D = "*[@]*.*";          E = $true  },
    @{ Row = 25; C = "foo@bar.com";        D = "?+[@]?+.?+";       E = $true  },
    @{ Row = 26; C = "foo@bar.com";        D = "@+[@]@+.@+";       E = $true  },
    @{ Row = 27; C = "+38(099)123-12-12";  D = "+7#(###)###-##-##"; E = $false },
    @{ Row = 28; C = "F";                  D = "F";                E = $true  },
    @{ Row = 29; C = "aBBBa";              D = "a*a";              E = $true  },
    @{ Row = 30; C = "F";                  D = "[A-Z]";            E = $true  },
    @{ Row = 31; C = "BAR+";               D = "[A-Z]++";          E = $true  },
    @{ Row = 32; C = "a2a";                D = "a#a";              E = $true  },
    @{ Row = 33; C = "aTa";                D = "a@a";              E = $true  },
    @{ Row = 34; C = "aM5b";               D = "a[L-P]#[!c-e]";    E = $true  },
    @{ Row = 35; C = "BAT123khg";          D = "B?T*";             E = $true  },
    @{ Row = 36; C = "AE1234AE";           D = "@@####@@";         E = $true  },
    @{ Row = 37; C = "123-45AE";           D = "###-##@@";         E = $true  },
    @{ Row = 38; C = "123-45AE";           D = "###-##@@";         E = $true  },
    @{ Row = 39; C = "123-45AE";           D = "###-##??+";        E = $true  },
    @{ Row = 40; C = "123-45AE123";        D = "###-##??+";        E = $true  },
    @{ Row = 41; C = "123-45-AE";          D = "#+-#+-@+";         E = $true  }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}

# Update the active selection recorded in the sheet view
$ws.Range("L30").Select()
